$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove a bunch of "spacer" empty paragraphs that were used between
#    sections of the MP 5319 body text. Two consecutive empty paragraphs
#    after "Mandatory Procedure" collapse to one; all the other isolated
#    empty paragraphs between sections are removed outright.
#    Deleting from the highest paragraph index down to the lowest keeps the
#    remaining indices stable while we work.
# ---------------------------------------------------------------------------
$emptyParaIndexesDescending = @(25, 23, 21, 19, 17, 15, 13, 11, 9, 7, 3)
foreach ($idx in $emptyParaIndexesDescending) {
    $d.Paragraphs.Item($idx).Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) Touch the "6. A copy of the solicitation's evaluation and award
#    factors." paragraph so the stale cached <w:lastRenderedPageBreak/>
#    marker (left over from before the paragraphs above were removed) gets
#    dropped when Word recalculates/re-saves the run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("6. A copy of the solicitation", $true, $false, $false, $false, $false, `
    $true, 1, $false, "6. A copy of the solicitation", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Style tweaks (spacing / keep-together / color) captured in styles.xml.
# ---------------------------------------------------------------------------

# Normal: add "space before" (12pt) alongside the existing "space after".
$normal = $d.Styles.Item("Normal")
$normal.ParagraphFormat.SpaceBefore = 12

# Heading 3: space-before grows from 12pt to 18pt (space-after goes away).
$heading3 = $d.Styles.Item("Heading 3")
$heading3.ParagraphFormat.SpaceBefore = 18
$heading3.ParagraphFormat.SpaceAfter = 0

# List 1: space before/after become symmetric (12pt/12pt) and contextual
# spacing (no space between paragraphs of the same style) is turned off.
$list1 = $d.Styles.Item("List 1")
$list1.ParagraphFormat.SpaceBefore = 12
$list1.ParagraphFormat.SpaceAfter = 12
$list1.NoSpaceBetweenParagraphsOfSameStyle = $false

# List 2 / List 3 / List 4: drop the explicit 6pt "space before" override
# (falls back to Normal's new 12pt/12pt) and turn contextual spacing off;
# keep-with-next / keep-together / the indent stay as they were.
foreach ($name in @("List 2", "List 3", "List 4")) {
    $s = $d.Styles.Item($name)
    $s.ParagraphFormat.SpaceBefore = 12
    $s.ParagraphFormat.SpaceAfter = 12
    $s.NoSpaceBetweenParagraphsOfSameStyle = $false
}

# List 1_change: keep-together turns on, the explicit spacing override goes
# away (falls back to Normal's new 12pt/12pt) and contextual spacing is
# turned off.
$list1change = $d.Styles.Item("List 1_change")
$list1change.ParagraphFormat.KeepTogether = $true
$list1change.ParagraphFormat.SpaceBefore = 12
$list1change.ParagraphFormat.SpaceAfter = 12
$list1change.NoSpaceBetweenParagraphsOfSameStyle = $false

# List 1_change Char: the red highlight color reverts to black.
$list1changeChar = $d.Styles.Item("List 1_change Char")
$list1changeChar.Font.Color = 0
